$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row by scanning column A (Beteckning).
$lastRow = 1
for ($r = 2; $r -le 5000; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    if ($a -eq $null -or $a -eq "") {
        break
    }
    $lastRow = $r
}

# Update column C ("Förändrad") date for every data row: 45184 -> 45186
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45186
    }
}

# For every row whose link columns (S, T, V, W, X, Y) hold a bare
# HYPERLINK(url) formula, add the record's "Beteckning" (column A) as the
# friendly display-text second argument: HYPERLINK(url, "A nnnnn-yyyy").
$cols = @("S", "T", "V", "W", "X", "Y")
for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$r")
        $f = $cell.Formula
        if ($f -like "*HYPERLINK(*" -and $f -notlike "*,*") {
            $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $name + '")'
            $cell.Formula = $newFormula
        }
    }
}
